# إضافة حدث جديد في Card19
# 1) Row 15 had several blank cells that should actually hold the literal
#    text "nan" (matching the rest of the sheet's "no data" convention).
# 2) A new service-event row (row 16) is appended for card 19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card19")

# --- Row 15: B..K and M were blank, they become literal text "nan" ---
# "nan" is not numeric-looking, so a plain assignment keeps it stored as text.
$ws.Range("B15:K15").Value = "nan"
$ws.Range("M15").Value = "nan"

# --- Row 16: new service event for card 19 ---
# "19" looks numeric, so force text with a leading apostrophe, then reset the
# style so no stray quote-prefix/number-format is left on the cell (the rest
# of the sheet carries no explicit cell style on data rows).
$ws.Range("A16").Value = "'19"
$ws.Range("A16").Style = "Normal"

# B16..K16 stay blank (present-but-empty cells, like the rest of the sheet).
# A lone apostrophe yields an empty *text* cell instead of no cell at all.
$ws.Range("B16:K16").Value = "'"
$ws.Range("B16:K16").Style = "Normal"

# L16: date-like text (backslash separators - Excel will not reparse it as a
# real date), plain assignment is safe here.
$ws.Range("L16").Value = "30\4\2025"

# M16: numeric-looking ("595.9"), force text the same way as A16.
$ws.Range("M16").Value = "'595.9"
$ws.Range("M16").Style = "Normal"

# N16 / O16: plain Arabic text, no numeric ambiguity.
$ws.Range("N16").Value = "تم تغيير الجرائد الاماميه ومعيارته (1_2_4_5_7_8) وسن السليندر"
$ws.Range("O16").Value = "الخبير"
